$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

# --- New row 9 data ---
$ws.Range("A9").Value = "Tree"
$ws.Range("B9").Value = "Hard"
$ws.Range("C9").Value = "124. Binary Tree Maximum Path Sum"

$note = "We need to find the max path, so there can be only 1 straight path in the left subtree and right subtree`nSo We'll do a post-order traversal because we need to find the maxPathSum for left and right subtree before we do calculations. And on each node of the subtrees, we have to add the max among the 3 choices, node.val itself, node.val + maxPathSum(left subtree) or node.val + maxPathSum(right subtree). Another Special case we need to consider is that the maxPath may not go through the root, it could go through any node, and thats what self.max_max will track, it'll track the max if the max path was passing through the current node from left to right (not up)"
$ws.Range("D9").Value = $note

# --- Formatting ---
# C9 gets the plain "Neutral" cell style (same family as C8 used, but without the extra
# vertical-center/wrap alignment that C8 has).
$ws.Range("C9").Style = "Neutral"

# D9 matches the other NOTES column cells (wrap text, top aligned).
$ws.Range("D9").WrapText = $true
$ws.Range("D9").VerticalAlignment = -4160

# C8's alignment loses the vertical-center (keeps wrap) now that the row below
# has its own plain variant of the Neutral style.
$ws.Range("C8").VerticalAlignment = -4107

# Row height for the new row (note column wraps to a taller row).
$ws.Rows.Item(9).RowHeight = 86.4

# --- Hyperlink for new problem link ---
$ws.Hyperlinks.Add($ws.Range("C9"), "https://leetcode.com/problems/binary-tree-maximum-path-sum/", "", "", "124. Binary Tree Maximum Path Sum") | Out-Null

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D9").Select()
